# Added new test cases for Change password from profile
# -----------------------------------------------------------------
# 1) Tiny tweak on the existing "changeName_happyPath" sheet
# 2) Two brand-new worksheets:
#      - checkChangePasswordValidation
#      - changePasswordHappyPath
# -----------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # changeName_happyPath
$ws2 = $wb.Worksheets.Item(2)   # profileChangeValidations

# --- sheet1: flip the runmode flag for the Azker/Moulana row to "N" ---
$ws1.Cells.Item(4,3).Value = "N"

# --- create the two new worksheets, appended after the last sheet ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws3.Name = "checkChangePasswordValidation"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws4.Name = "changePasswordHappyPath"

# ===================================================================
# Populate ws3 ("checkChangePasswordValidation") — the exact order of
# these writes matters: it controls the order new shared strings are
# minted in, so keep it as-is.
# ===================================================================
$ws3.Cells.Item(1,1).Value = "currentpwd"
$ws3.Cells.Item(1,2).Value = "newpwd"
$ws3.Cells.Item(1,3).Value = "confirmpwd"
$ws3.Cells.Item(2,1).Value = "abc"
$ws3.Cells.Item(2,2).Value = "xyz"
$ws3.Cells.Item(1,4).Value = "message"
$ws3.Cells.Item(2,4).Value = "Current password is incorrect"
$ws3.Cells.Item(5,2).Value = "Abcd@1234"
$ws3.Cells.Item(3,4).Value = "Your password must contain minimum 8 characters"
$ws3.Cells.Item(3,1).Value = "Intel@123"
$ws3.Cells.Item(4,2).Value = "abcdefgh"
$ws3.Cells.Item(4,4).Value = "Your password must contain uppercase and numeric characters"
$ws3.Cells.Item(5,4).Value = "Passwords you provided doesn't match"
$ws3.Cells.Item(5,3).Value = "Abcd@123"

# remaining cells reuse already-minted shared strings
$ws3.Cells.Item(1,5).Value = "runmode"
$ws3.Cells.Item(2,3).Value = "xyz"
$ws3.Cells.Item(2,5).Value = "N"
$ws3.Cells.Item(3,2).Value = "abc"
$ws3.Cells.Item(3,3).Value = "abc"
$ws3.Cells.Item(3,5).Value = "N"
$ws3.Cells.Item(4,1).Value = "Intel@123"
$ws3.Cells.Item(4,3).Value = "abcdefgh"
$ws3.Cells.Item(4,5).Value = "N"
$ws3.Cells.Item(5,1).Value = "Intel@123"
$ws3.Cells.Item(5,5).Value = "Y"

# hyperlink-styled cells on ws3
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), "https://example.com/checkChangePasswordValidation/row3") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), "https://example.com/checkChangePasswordValidation/row4") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), "https://example.com/checkChangePasswordValidation/row5a") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,2), "https://example.com/checkChangePasswordValidation/row5b") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,3), "https://example.com/checkChangePasswordValidation/row5c") | Out-Null

$ws3.Columns.Item(1).ColumnWidth = 9.983072916666666
$ws3.Columns.Item(2).ColumnWidth = 18.436197916666668
$ws3.Columns.Item(3).ColumnWidth = 15.256510416666666
$ws3.Columns.Item(4).ColumnWidth = 44.346354166666664

$ws3.Range("A1:E3").Select() | Out-Null

# ===================================================================
# Populate ws4 ("changePasswordHappyPath")
# ===================================================================
$ws4.Cells.Item(1,1).Value = "currentpwd"
$ws4.Cells.Item(1,2).Value = "newpwd"
$ws4.Cells.Item(1,3).Value = "confirmpwd"
$ws4.Cells.Item(1,4).Value = "message"
$ws4.Cells.Item(2,4).Value = "Password changed!You successfully changed your password."

$ws4.Cells.Item(1,5).Value = "runmode"
$ws4.Cells.Item(2,1).Value = "Intel@123"
$ws4.Cells.Item(2,2).Value = "Intel@123"
$ws4.Cells.Item(2,3).Value = "Intel@123"
$ws4.Cells.Item(2,5).Value = "Y"

$ws4.Hyperlinks.Add($ws4.Cells.Item(2,1), "https://example.com/changePasswordHappyPath/row2a") | Out-Null
$ws4.Hyperlinks.Add($ws4.Cells.Item(2,2), "https://example.com/changePasswordHappyPath/row2b") | Out-Null
$ws4.Hyperlinks.Add($ws4.Cells.Item(2,3), "https://example.com/changePasswordHappyPath/row2c") | Out-Null

$ws4.Columns.Item(1).ColumnWidth = 13.072916666666666
$ws4.Columns.Item(2).ColumnWidth = 8.799479166666666
$ws4.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws4.Columns.Item(4).ColumnWidth = 52.529947916666664

# Final selections: sheet1 -> D12, then land on ws4 (last active tab == ws4)
$ws1.Range("D12").Select() | Out-Null
$ws2.Range("B9").Select() | Out-Null
$ws4.Cells.Item(16,8).Select() | Out-Null

Write-Output "edit applied"
